# Improve logging system configuration.
# Appends a new data row (row 60) recorded at 2025-07-08 11:47:00 to each
# of the four log sheets, mirroring the existing row layout/format.

$wb = $excel.ActiveWorkbook

$dateValue = 45846.49097222222
$dateFormat = "YYYY-MM-DD HH:MM:SS"

$rows = @(
    @{
        Sheet = 1
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x50"
        E = "0xf"
        F = 380
        G = "7.598631275147109e+23"
        H = 336
        I = 15
    },
    @{
        Sheet = 2
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x60"
        E = "0xe"
        F = 400
        G = "5.68432987514711e+23"
        H = 352
        I = 14
    },
    @{
        Sheet = 3
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x67"
        E = "0x3"
        F = 110
        G = "5.68631262647114e+23"
        H = 103
        I = 3
    },
    @{
        Sheet = 4
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x67"
        E = "0x3"
        F = 110
        G = "9.85046333984776e+23"
        H = 103
        I = 3
    }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)
    $newRow = 60

    $ws.Cells.Item($newRow, 1).Value = $dateValue
    $ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($newRow, 2).Value = $row.B
    $ws.Cells.Item($newRow, 3).Value = $row.C
    $ws.Cells.Item($newRow, 4).Value = $row.D
    $ws.Cells.Item($newRow, 5).Value = $row.E
    $ws.Cells.Item($newRow, 6).Value = $row.F
    $ws.Cells.Item($newRow, 7).Value = [double]$row.G
    $ws.Cells.Item($newRow, 8).Value = $row.H
    $ws.Cells.Item($newRow, 9).Value = $row.I
}
